# Applies the commit's changes to "Especificaciones - Copia.xlsx":
#  - Updates a couple of "Estado"/"Comentarios" cells on the
#    "Cálculos Sectores" sheet with new progress notes.
#  - Leaves the sheet scrolled/selected near the rows that were edited,
#    mirroring the view state the author left the workbook in when saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cálculos Sectores")
$ws.Activate()

# Row 9 (ID 3080): Estado -> "Creo que algo falla.", Comentarios -> balance/result note.
$ws.Range("C9").Value = "Creo que algo falla."
$ws.Range("D9").Value = "El balance de error de los sectores y el resultado de los calculos son algo altos…"

# Row 11 (ID 3090, second entry): Estado -> "casi terminado",
# Comentarios -> question about needing an ideal power formula.
$ws.Range("C11").Value = "casi terminado"
$ws.Range("D11").Value = "ERROR DE POTENCIA? Para eso seria necesario una potencia ideal, según la formula que tengo."

# Leave the window scrolled/selected where the author left it after editing:
# selection on C12, viewport scrolled up/right so row 9 / column B is the
# top-left visible cell.
$ws.Range("C12").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 2
